$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New journal entries (rows 19 and 20) ---
$ws.Range("B19").Value = "MA-20"
$ws.Range("C19").Value = [datetime]"2019-04-04"
$ws.Range("D19").Value = "90m"
$ws.Range("E19").Value = "Codage"
$ws.Range("F19").Value = "Finalisation du code"

$ws.Range("B20").Value = "MA-20"
$ws.Range("C20").Value = [datetime]"2019-04-05"
$ws.Range("D20").Value = "90m"
$ws.Range("E20").Value = "Codage"
$ws.Range("F20").Value = "Finalisation du code"

# --- Update the view: scroll position + active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("H20").Select()
